$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Nota" column (F). This shifts the old "Texto" column (G) into F.
$ws.Range("F1:F2").Delete()

# Update row 2 values to the new data
$ws.Range("B2").Value = "Bom Dia Inter"
$ws.Range("C2").Value = "PROCON"
$ws.Range("D2").Value = "2025-03-31T17:41"
$ws.Range("E2").Value = "Neutro"
$ws.Range("F2").Value = "aaateste"
